# Add a new worksheet "ODI Bowling Extra" (mirrors the existing
# "ODI Batting Extra" sheet, but for bowling-related scraped attributes)
# and drop the stray empty E8 cell on "ODI Batting Extra".

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "ODI Batting Extra": remove the empty E8 cell
# ------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("E8").ClearContents()

# ------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" worksheet after the existing
#    last sheet ("ODI Batting Extra")
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Store everything as text (matches the scraped source data, which
# keeps match codes / percentages / counts as plain strings).
$ws.Range("A1:C21").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "MAIDEN_OVERS"
$ws.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# Data rows
$rows = @(
    @("4572", "0", "10.00%"),
    @("4573", "", ""),
    @("4575", "0", "10.00%"),
    @("4576", "0", "20.00%"),
    @("4578", "0", ""),
    @("4581", "1", ""),
    @("4604", "0", ""),
    @("4625", "0", "20.00%"),
    @("4629", "1", "20.00%"),
    @("4631", "", ""),
    @("4632", "0", "10.00%"),
    @("4635", "", ""),
    @("4677", "0", "10.00%"),
    @("4681", "0", "10.00%"),
    @("4680", "0", ""),
    @("4684", "0", ""),
    @("4702", "1", "40.00%"),
    @("4703", "0", "10.00%"),
    @("4705", "1", "30.00%"),
    @("4706", "", "")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("A1").Select()
